$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C16").Value = "vinculo-pacientes"
$ws.Range("G7").Value = "pacientes-disciplina"
$ws.Range("G8").Value = "ano/semestre"

$ws.Range("G9").Select()
